$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Project 4 - Report Task " + "4" (two runs) -> one run with the
#    combined text "Project 4 - Report Task 4".
#    A plain Find/Replace over the whole phrase merges the (identically
#    formatted) runs it touches into a single run, which is exactly the
#    shape the target XML wants here.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Project 4 " + [char]0x2013 + " Report Task 4", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Project 4 " + [char]0x2013 + " Report Task 4", 2)

# ------------------------------------------------------------------
# 2) "Alex (Sean) Wall" -> "Sean (Alex) Wall", but split into four
#    separate (identically formatted) runs: "Sean", " (", "Alex", ") Wall".
#    A simple Find/Replace would merge same-formatted text back into a
#    single run, so instead we rebuild the phrase using Copy/Paste of
#    the original sub-ranges (Word keeps pasted fragments as their own
#    runs even when formatting matches the surroundings), then delete
#    the original leftover text.
# ------------------------------------------------------------------
$p = $d.Paragraphs(3)
$r = $p.Range

$full = $d.Range($r.Start, $r.End)
$full.Find.Execute("Alex (Sean) Wall", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$origStart = $full.Start

# Offsets within the original "Alex (Sean) Wall" text:
#   Alex = [0,4)   " (" = [4,6)   Sean = [6,10)   ") Wall" = [10,16)
$shift = 0

# a) "Sean"
$src = $d.Range($origStart + $shift + 6, $origStart + $shift + 10)
$src.Copy()
$ins = $d.Range($origStart + $shift, $origStart + $shift)
$ins.Paste()
$shift = $shift + 4

# b) " ("
$src = $d.Range($origStart + $shift + 4, $origStart + $shift + 6)
$src.Copy()
$ins = $d.Range($origStart + $shift, $origStart + $shift)
$ins.Paste()
$shift = $shift + 2

# c) "Alex"
$src = $d.Range($origStart + $shift + 0, $origStart + $shift + 4)
$src.Copy()
$ins = $d.Range($origStart + $shift, $origStart + $shift)
$ins.Paste()
$shift = $shift + 4

# d) ") Wall"
$src = $d.Range($origStart + $shift + 10, $origStart + $shift + 16)
$src.Copy()
$ins = $d.Range($origStart + $shift, $origStart + $shift)
$ins.Paste()
$shift = $shift + 6

# Remove the now-redundant original "Alex (Sean) Wall" text that got
# pushed to the right of everything we just pasted.
$leftover = $d.Range($origStart + $shift, $origStart + $shift + 16)
$leftover.Delete()
